# quarterly.xlsx update: roll the quarter window forward one quarter.
#  - drop the oldest quarter column (1399/06)
#  - shift all remaining quarters one column to the left (E..N)
#  - add the newest quarter (1401/12) in column N
#  - apply the restated figures for 1400/12 that came with the new
#    read_price algorithm (rows 15, 19, 26, 27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N")

# Header rows (row 8 and row 24) hold the quarter labels.
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $quarters[$i]
    $ws.Range($cols[$i] + "24").Value = $quarters[$i]
}

# Data rows: new value set for each row, columns E..N (one quarter shift + new
# quarter in N; rows 15/19/26/27 also carry a restated 1400/12 figure in J).
$rowData = @{
    10 = @(92,230,0,0,352,-115,0,0,200,-200)
    15 = @(45,162,48,77,67,150,86,166,127,113)
    16 = @(750,810,881,888,897,3182,1022,2938,2240,2074)
    17 = @(7151,9838,9627,13529,12104,86670,15785,71971,44223,54156)
    19 = @(8150,29657,18535,39408,17044,-25394,29878,3719,1528,35604)
    20 = @(16188,40697,29091,53902,30464,64493,46771,78794,48318,91747)
    26 = @(97,95,50,95,95,93,89,95,95,93)
    27 = @(301,325,367,329,329,342,300,307,305,334)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
